$d = $word.ActiveDocument

# The "Requisitos" bullet list currently reads (in order):
#   LOB1006 -  Cálculo IV  (Requisito fraco)
#   LOB1019 -  Física II  (Requisito fraco)
#   LOQ4083 -  Fenômenos de Transporte I  (Requisito fraco)
# It needs to become:
#   LOQ4083 -  Fenômenos de Transporte I  (Requisito fraco)
#   LOB1006 -  Cálculo IV  (Requisito fraco)
#   LOB1019 -  Física II  (Requisito fraco)
# i.e. move the LOQ4083 line (its own run, text + manual line break) from the
# end of the list to the front, leaving the other two runs untouched.

# 1) Remove the existing LOQ4083 line (text + its manual line break) from the
#    end of the requirements list.
$old = $d.Content
$foundOld = $old.Find.Execute("LOQ4083 -  Fenômenos de Transporte I  (Requisito fraco)^l", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundOld) {
    $old.Delete()
}

# 2) Re-insert it as a standalone run immediately before the LOB1006 line,
#    using InsertXML so it lands as its own <w:r> (text + <w:br/>) rather
#    than merging into the neighbouring run.
$dest = $d.Content
$foundDest = $dest.Find.Execute("LOB1006 -  Cálculo IV  (Requisito fraco)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundDest) {
    $insertPoint = $d.Range($dest.Start, $dest.Start)
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>LOQ4083 -  Fenômenos de Transporte I  (Requisito fraco)</w:t><w:br/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $insertPoint.InsertXML($xml)
}
